$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data sources in column D ("Link") are updated from plain URLs to
# HTML-style anchor tags, e.g. <a href='...'>ONS</a>
$ws.Range("D3").Value = "<a href='https://www.ons.gov.uk/peoplepopulationandcommunity/wellbeing/articles/subnationalindicatorsexplorer/2022-01-06'>ONS</a>"
$ws.Range("D4").Value = "<a href='https://www.aoc.co.uk/research-unit/data-sources'>AOC</a>"
$ws.Range("D5").Value = "<a href='https://census.gov.uk/local-authorities'>Census</a>"
$ws.Range("D6").Value = "<a href='https://www.nomisweb.co.uk/'>Nomis</a>"
$ws.Range("D7").Value = "<a href='https://explore-education-statistics.service.gov.uk/'>EES</a>"

# The selected cell moved from B14 to E10 when the file was last saved
$ws.Range("E10").Select()
